$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @(0, 0, 0, 0.442943416743802, 0, 0.068008934099110505, 0, 0, 0, 0.40641146779420501, 0, 0.0826361813628822),
    @(0.52667637912673504, 0, 0.0293348611438717, 0, 0, 0, 0.40014397938859902, 0, 0.043844780340794802, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0.40989169877151099, 0, 0, 0, 0, 0, 0.59010830122848901),
    @(0, 0, 0.34477421302384897, 0, 0, 0, 0, 0, 0.65522578697615097, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0.656774168989399, 0, 0.343225831010601),
    @(0, 0, 0, 0.80745471387803402, 0, 0.192545286121966, 0, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 0.72641020391599898, 0, 0.27358979608400102, 0, 0, 0),
    @(0.76759020332927297, 0, 0.232409796670727, 0, 0, 0, 0, 0, 0, 0, 0, 0)
)

for ($r = 1; $r -le 8; $r++) {
    for ($c = 1; $c -le 12; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$r - 1][$c - 1]
    }
}

$ws.Range("H20").Select()
